# Commit message: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted at row 63 of the sheet,
# pushing all existing data rows (previously 63-119) down by one
# (now 64-120). The sheet's used range grows from A1:R119 to A1:R120.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 63; everything below (old rows
# 63-119) shifts down to 64-120, carrying its existing formatting.
$ws.Rows(63).Insert()

# Populate the newly inserted row 63 with the new observation.
$ws.Range("A63").Value = 4
$ws.Range("B63").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C63").Value = "Los Lagos"
$ws.Range("D63").Value = 44740
$ws.Range("E63").Value = 10
$ws.Range("F63").Value = 100112022
$ws.Range("G63").Value = "Arveja Verde"
$ws.Range("H63").Value = "Perfection"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 40
$ws.Range("K63").Value = 44000
$ws.Range("L63").Value = 44000
$ws.Range("M63").Value = 44000
$ws.Range("N63").Value = "`$/malla 25 kilos"
$ws.Range("O63").Value = "Provincia de Limarí"
$ws.Range("P63").Value = 1760
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
